$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Semestre ideal: EQD-5,EQN-6 -> EQD-8,EQN-9
$ws.Range("B9").Value = "EQD-8,EQN-9"
$ws.Range("C9").Value = "EQD-8,EQN-9"

# 2) Add a new requirement row above the existing one. Push the current
#    requirement row (row 24: "LOQ4055 ... (Requisito fraco)") down to row
#    25, carrying its formatting and row height with it...
$ws.Range("B24:C24").Copy($ws.Range("B25:C25"))
$ws.Rows.Item(25).RowHeight = $ws.Rows.Item(24).RowHeight()

# ...then put the new requirement text ("LOQ4002 ... (Requisito fraco)") into
# row 24, which already has the right style/height for a requirement row.
$ws.Range("B24").Value = "LOQ4002 -  Reatores Quimicos  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOQ4002 -  Reatores Quimicos  (Requisito fraco)`n"

Write-Output "done"
